# Actualización automática 2025-06-24 15:25:08
#
# A new advisor, "AREVALO PEÑA JORGE LUIS", is inserted as the first data
# row (row 2) on both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets.
# All existing advisor rows shift down by one row, and the trailing summary
# row (the "X de N" counts row on VENTAS POR GRUPO, and the totals row on
# VENTA MENSUAL) moves down with them. The summary row's "de 5" counts
# become "de 6" to reflect the now 6-advisor roster.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row above the current row 2 (AVILA TORRES ...), pushing
# everything (including the trailing "X de 5" row) down by one.
$ws1.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above (the bold
# header); strip that back to the plain data-row look used by the other
# advisor rows, then restore the currency number format on the numeric
# columns so it reuses the workbook's existing currency style.
$ws1.Rows.Item(2).ClearFormats()
$ws1.Range("C2:R2").NumberFormat = """$""#,##0.00"

$ws1.Range("A2").Value = "OFICINA-CATAECSA"
$ws1.Range("B2").Value = "AREVALO PEÑA JORGE LUIS"
$ws1.Range("C2:R2").Value = 0

# The old "0 de 5" / "1 de 5" summary row is now row 8; bump the counts
# to reflect the 6th advisor that was just added.
$ws1.Range("C8").Value = "0 de 6"
$ws1.Range("D8").Value = "1 de 6"
$ws1.Range("E8").Value = "0 de 6"
$ws1.Range("F8").Value = "0 de 6"
$ws1.Range("G8").Value = "0 de 6"
$ws1.Range("H8").Value = "0 de 6"
$ws1.Range("I8").Value = "0 de 6"
$ws1.Range("J8").Value = "0 de 6"
$ws1.Range("K8").Value = "0 de 6"
$ws1.Range("L8").Value = "0 de 6"
$ws1.Range("M8").Value = "0 de 6"
$ws1.Range("N8").Value = "0 de 6"
$ws1.Range("O8").Value = "0 de 6"
$ws1.Range("P8").Value = "0 de 6"
$ws1.Range("Q8").Value = "0 de 6"
$ws1.Range("R8").Value = "0 de 6"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(2).Insert()
$ws2.Rows.Item(2).ClearFormats()
$ws2.Range("C2:G2").NumberFormat = """$""#,##0.00"

$ws2.Range("A2").Value = "OFICINA-CATAECSA"
$ws2.Range("B2").Value = "AREVALO PEÑA JORGE LUIS"
$ws2.Range("C2:G2").Value = 0
